$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "52÷8=6, 4"
$t.Cell(1,2).Range.Text = "79÷2=39, 1"
$t.Cell(1,3).Range.Text = "69÷5=13, 4"
$t.Cell(1,4).Range.Text = "77÷4=19, 1"
$t.Cell(1,5).Range.Text = "52÷4=13, 0"
$t.Cell(5,1).Range.Text = "39÷2=19, 1"
$t.Cell(5,2).Range.Text = "92÷9=10, 2"
$t.Cell(5,3).Range.Text = "69÷2=34, 1"
$t.Cell(5,4).Range.Text = "40÷3=13, 1"
$t.Cell(5,5).Range.Text = "70÷7=10, 0"
$t.Cell(9,1).Range.Text = "27÷9=3, 0"
$t.Cell(9,2).Range.Text = "52÷7=7, 3"
$t.Cell(9,3).Range.Text = "23÷2=11, 1"
$t.Cell(9,4).Range.Text = "44÷4=11, 0"
$t.Cell(9,5).Range.Text = "21÷3=7, 0"
$t.Cell(13,1).Range.Text = "34÷8=4, 2"
$t.Cell(13,2).Range.Text = "48÷3=16, 0"
$t.Cell(13,3).Range.Text = "27÷6=4, 3"
$t.Cell(13,4).Range.Text = "79÷9=8, 7"
$t.Cell(13,5).Range.Text = "62÷7=8, 6"
$t.Cell(17,1).Range.Text = "55÷8=6, 7"
$t.Cell(17,2).Range.Text = "14÷6=2, 2"
$t.Cell(17,3).Range.Text = "46÷4=11, 2"
$t.Cell(17,4).Range.Text = "63÷9=7, 0"
$t.Cell(17,5).Range.Text = "94÷7=13, 3"
